$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-04 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-05 Wednesday", 2)

$d.Content.Find.Execute("126×7=882", $true, $false, $false, $false, $false, $true, 1, $false, "633×4=2532", 2)
$d.Content.Find.Execute("889×7=6223", $true, $false, $false, $false, $false, $true, 1, $false, "476×5=2380", 2)
$d.Content.Find.Execute("360×9=3240", $true, $false, $false, $false, $false, $true, 1, $false, "640×8=5120", 2)
$d.Content.Find.Execute("602×7=4214", $true, $false, $false, $false, $false, $true, 1, $false, "887×7=6209", 2)
$d.Content.Find.Execute("969×4=3876", $true, $false, $false, $false, $false, $true, 1, $false, "266×9=2394", 2)
$d.Content.Find.Execute("567×5=2835", $true, $false, $false, $false, $false, $true, 1, $false, "846×7=5922", 2)
$d.Content.Find.Execute("685×5=3425", $true, $false, $false, $false, $false, $true, 1, $false, "305×9=2745", 2)
$d.Content.Find.Execute("177×9=1593", $true, $false, $false, $false, $false, $true, 1, $false, "781×3=2343", 2)
$d.Content.Find.Execute("406×3=1218", $true, $false, $false, $false, $false, $true, 1, $false, "666×6=3996", 2)
$d.Content.Find.Execute("702×9=6318", $true, $false, $false, $false, $false, $true, 1, $false, "490×9=4410", 2)
$d.Content.Find.Execute("531×4=2124", $true, $false, $false, $false, $false, $true, 1, $false, "880×3=2640", 2)
$d.Content.Find.Execute("991×2=1982", $true, $false, $false, $false, $false, $true, 1, $false, "552×7=3864", 2)
$d.Content.Find.Execute("241×4=964", $true, $false, $false, $false, $false, $true, 1, $false, "193×2=386", 2)
$d.Content.Find.Execute("735×3=2205", $true, $false, $false, $false, $false, $true, 1, $false, "499×3=1497", 2)
$d.Content.Find.Execute("585×5=2925", $true, $false, $false, $false, $false, $true, 1, $false, "368×3=1104", 2)
$d.Content.Find.Execute("357×9=3213", $true, $false, $false, $false, $false, $true, 1, $false, "346×6=2076", 2)
$d.Content.Find.Execute("451×7=3157", $true, $false, $false, $false, $false, $true, 1, $false, "139×4=556", 2)
$d.Content.Find.Execute("832×9=7488", $true, $false, $false, $false, $false, $true, 1, $false, "233×7=1631", 2)
$d.Content.Find.Execute("525×6=3150", $true, $false, $false, $false, $false, $true, 1, $false, "690×8=5520", 2)
$d.Content.Find.Execute("583×3=1749", $true, $false, $false, $false, $false, $true, 1, $false, "816×7=5712", 2)
$d.Content.Find.Execute("338×2=676", $true, $false, $false, $false, $false, $true, 1, $false, "788×5=3940", 2)
$d.Content.Find.Execute("151×9=1359", $true, $false, $false, $false, $false, $true, 1, $false, "972×7=6804", 2)
$d.Content.Find.Execute("176×5=880", $true, $false, $false, $false, $false, $true, 1, $false, "321×6=1926", 2)
$d.Content.Find.Execute("912×4=3648", $true, $false, $false, $false, $false, $true, 1, $false, "364×8=2912", 2)
$d.Content.Find.Execute("483×2=966", $true, $false, $false, $false, $false, $true, 1, $false, "384×3=1152", 2)
